$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameter")

# Update the comment text on A2 (the jx:each directive) so the query
# filters/binds on customer number instead of customer name.
# Single-quoted here-string avoids any PowerShell $-interpolation of
# the literal ${cusNumber} token.
$newCommentText = @'
jx:each(items="jdbc.query('SELECT CUS_CUSTOMERNAME name, CUS_PHONE phone FROM T_AGG_CUSTOMER WHERE CUS_CUSTOMERNUMBER = ${cusNumber}')" var="customer" lastCell="B2")

'@
$cmt = $ws.Range("A2").Comment
$cmt.Text($newCommentText)

# Update the parameter help text cells
$ws.Range("A4").Value = '(You need a text parameter in ReportServer with key=cusNumber. You can then pass any customer name, e.g. 350)'
$ws.Range("A5").Value = '(Version: 1.0.1)'
$ws.Range("A6").Value = '(Last tested with: ReportServer 4.0.0-6053) '

# Update the selected/active cell in the sheet view
$ws.Range("A5").Select()
